$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 399, shifting existing rows 399-417 down to 400-418
$ws.Rows.Item(399).Insert()

# Populate the newly inserted row 399 with the new data record
$ws.Cells.Item(399, 1).Value = 9
$ws.Cells.Item(399, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(399, 3).Value = "Metropolitana"
$ws.Cells.Item(399, 4).Value = 44747
$ws.Cells.Item(399, 4).Style = $ws.Cells.Item(400, 4).Style
$ws.Cells.Item(399, 4).NumberFormat = $ws.Cells.Item(400, 4).NumberFormat
$ws.Cells.Item(399, 5).Value = 13
$ws.Cells.Item(399, 6).Value = 100112039
$ws.Cells.Item(399, 7).Value = "Ciboulette"
$ws.Cells.Item(399, 8).Value = "Sin especificar"
$ws.Cells.Item(399, 9).Value = "Primera"
$ws.Cells.Item(399, 10).Value = 250
$ws.Cells.Item(399, 11).Value = 900
$ws.Cells.Item(399, 12).Value = 1000
$ws.Cells.Item(399, 13).Value = 950
$ws.Cells.Item(399, 14).Value = "$/docena de atados"
$ws.Cells.Item(399, 15).Value = "Región Metropolitana"
$ws.Cells.Item(399, 16).Value = 317
$ws.Cells.Item(399, 17).Value = 3
$ws.Cells.Item(399, 18).Value = "Hortaliza"
